$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.301.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").Value = "'3.898.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'485.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").Value = "'145.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("D7").Value = "'0.623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = "'0.997"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +2.59%  '
$ws.Range("D10").Value = "'0.178"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.70%  '
$ws.Range("D11").Value = "'0.0000354"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.67%  '
$ws.Range("D12").Value = "'42.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").Value = "'4.519.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").Value = "'3.922.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("D16").Value = "'14.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("D18").Value = "'19.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D20").Value = "'68.313.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.12%  '
$ws.Range("D21").Value = "'430.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.87%  '
$ws.Range("E22").Value = '  +7.10%  '
$ws.Range("E23").Value = '  +0.99%  '
$ws.Range("D24").Value = "'12.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +20.03%  '
$ws.Range("D25").Value = "'88.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("E26").Value = '  +4.33%  '
$ws.Range("E27").Value = '  -5.37%  '
$ws.Range("D28").Value = "'37.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.55%  '
$ws.Range("E29").Value = '  -3.62%  '
$ws.Range("D30").Value = "'719.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'13.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.82%  '
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("E33").Value = '  +2.66%  '
$ws.Range("D34").Value = "'61.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.17%  '
$ws.Range("D35").Value = "'0.0₃0874"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.38%  '
$ws.Range("D36").Value = "'6.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.01%  '
$ws.Range("D37").Value = "'40.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("D38").Value = "'0.398"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +17.49%  '
$ws.Range("E39").Value = '  -3.83%  '
$ws.Range("D40").Value = "'0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = "'0.0498"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.83%  '
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").Value = "'2.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.52%  '
$ws.Range("D43").Value = "'3.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.08%  '
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = "'0.0₆0370"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +33.22%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = "'0.142"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("E48").Value = '  +5.60%  '
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("E50").Value = '  -2.39%  '
$ws.Range("D51").Value = "'144.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.50%  '